$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H92").Value = 634.0833
$ws_ALC.Range("I92").Value = 634
$ws_ALC.Range("J92").Value = 634.5
$ws_ALC.Range("K92").Value = 634
$ws_ALC.Range("L92").Value = 634.5
$ws_ALC.Range("M92").Value = 614
$ws_ALC.Range("N92").Value = -3130.5
$ws_ALC.Range("H97").Value = 1804
$ws_ALC.Range("I97").Value = 0
$ws_ALC.Range("K97").Value = 0
$ws_ALC.Range("M97").ClearContents()
$ws_ALC.Range("H137").Value = 64648
$ws_ALC.Range("I137").Value = 2364.1667
$ws_ALC.Range("J137").Value = 251499.5
$ws_ALC.Range("K137").Value = 7092.500100000001
$ws_ALC.Range("L137").Value = 754498.5
$ws_ALC.Range("M137").Value = -4542.500100000001
$ws_ALC.Range("N137").Value = -759598.5
$ws_ALC.Range("H138").Value = 1622.7222
$ws_ALC.Range("I138").Value = 1156.6897
$ws_ALC.Range("K138").Value = 3470.0691
$ws_ALC.Range("M138").Value = 1669.9309
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 32893.926
$ws_ARM.Range("I32").Value = 19329.527
$ws_ARM.Range("J32").Value = 100715.91
$ws_ARM.Range("K32").Value = 19329.527
$ws_ARM.Range("L32").Value = 100715.91
$ws_ARM.Range("M32").Value = -19042.527
$ws_ARM.Range("N32").Value = -101289.91
$ws_ARM.Range("H45").Value = 844668.5
$ws_ARM.Range("I45").Value = 2022403.8
$ws_ARM.Range("J45").Value = 3429
$ws_ARM.Range("K45").Value = 2022403.8
$ws_ARM.Range("L45").Value = 3429
$ws_ARM.Range("M45").Value = -2022026.8
$ws_ARM.Range("N45").Value = -4183
$ws_ARM.Range("H61").Value = 3046.5833
$ws_ARM.Range("I61").Value = 2444.875
$ws_ARM.Range("K61").Value = 2444.875
$ws_ARM.Range("M61").Value = -2232.875
$ws_ARM.Range("H132").Value = 11608.429
$ws_ARM.Range("J132").Value = 5399.8
$ws_ARM.Range("L132").Value = 16199.4
$ws_ARM.Range("N132").Value = -21259.4
$ws_ARM.Range("H136").Value = 3046.5833
$ws_ARM.Range("I136").Value = 2444.875
$ws_ARM.Range("K136").Value = 7334.625
$ws_ARM.Range("M136").Value = -4784.625
$ws_ARM.Range("H141").Value = 98462.336
$ws_ARM.Range("J141").Value = 98462.336
$ws_ARM.Range("L141").Value = 98462.336
$ws_ARM.Range("N141").Value = -108822.336
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H20").Value = 5006657.5
$ws_BSM.Range("I20").Value = 10424746
$ws_BSM.Range("J20").Value = 5345.231
$ws_BSM.Range("K20").Value = 10424746
$ws_BSM.Range("L20").Value = 5345.231
$ws_BSM.Range("M20").Value = -10424499
$ws_BSM.Range("N20").Value = -5839.231
$ws_BSM.Range("H35").Value = 30000
$ws_BSM.Range("J35").Value = 30000
$ws_BSM.Range("L35").Value = 30000
$ws_BSM.Range("N35").Value = -30620
$ws_BSM.Range("H86").Value = 2216.9487
$ws_BSM.Range("I86").Value = 1794.7391
$ws_BSM.Range("K86").Value = 1794.7391
$ws_BSM.Range("M86").Value = -671.7391
$ws_BSM.Range("H89").Value = 2216.9487
$ws_BSM.Range("I89").Value = 1794.7391
$ws_BSM.Range("K89").Value = 8973.6955
$ws_BSM.Range("M89").Value = -3357.6955
$ws_BSM.Range("H105").Value = 3573657.5
$ws_BSM.Range("I105").Value = 4002136.2
$ws_BSM.Range("K105").Value = 4002136.2
$ws_BSM.Range("M105").Value = -4000389.2
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H9").Value = 91040.336
$ws_CRP.Range("J9").Value = 91040.336
$ws_CRP.Range("L9").Value = 91040.336
$ws_CRP.Range("N9").Value = -91376.336
$ws_CRP.Range("H16").Value = 3092
$ws_CRP.Range("I16").Value = 3092
$ws_CRP.Range("J16").Value = 0
$ws_CRP.Range("K16").Value = 3092
$ws_CRP.Range("L16").Value = 0
$ws_CRP.Range("M16").Value = -2805
$ws_CRP.Range("N16").ClearContents()
$ws_CRP.Range("H31").Value = 2193
$ws_CRP.Range("I31").Value = 1927.1904
$ws_CRP.Range("K31").Value = 1927.1904
$ws_CRP.Range("M31").Value = -1632.1904
$ws_CRP.Range("H34").Value = 2193
$ws_CRP.Range("I34").Value = 1927.1904
$ws_CRP.Range("K34").Value = 1927.1904
$ws_CRP.Range("M34").Value = -1725.1904
$ws_CRP.Range("H41").Value = 12666.167
$ws_CRP.Range("J41").Value = 12999.25
$ws_CRP.Range("L41").Value = 12999.25
$ws_CRP.Range("N41").Value = -13855.25
$ws_CRP.Range("H58").Value = 5862.391
$ws_CRP.Range("I58").Value = 8887
$ws_CRP.Range("K58").Value = 8887
$ws_CRP.Range("M58").Value = -8684
$ws_CRP.Range("H113").Value = 3092
$ws_CRP.Range("I113").Value = 3092
$ws_CRP.Range("J113").Value = 0
$ws_CRP.Range("K113").Value = 3092
$ws_CRP.Range("L113").Value = 0
$ws_CRP.Range("M113").Value = -922
$ws_CRP.Range("N113").ClearContents()
$ws_CRP.Range("H122").Value = 5049.9165
$ws_CRP.Range("I122").Value = 5257
$ws_CRP.Range("J122").Value = 4760
$ws_CRP.Range("K122").Value = 15771
$ws_CRP.Range("L122").Value = 14280
$ws_CRP.Range("M122").Value = -13321
$ws_CRP.Range("N122").Value = -19180
$ws_CRP.Range("H132").Value = 2934
$ws_CRP.Range("I132").Value = 2729.65
$ws_CRP.Range("K132").Value = 8188.950000000001
$ws_CRP.Range("M132").Value = -5658.950000000001
$ws_CRP.Range("H133").Value = 50000
$ws_CRP.Range("J133").Value = 50000
$ws_CRP.Range("L133").Value = 50000
$ws_CRP.Range("N133").Value = -55060
$ws_CRP.Range("H134").Value = 3876.1904
$ws_CRP.Range("I134").Value = 3689.4707
$ws_CRP.Range("K134").Value = 11068.4121
$ws_CRP.Range("M134").Value = -8533.4121
$ws_CRP.Range("H136").Value = 5862.391
$ws_CRP.Range("I136").Value = 8887
$ws_CRP.Range("K136").Value = 26661
$ws_CRP.Range("M136").Value = -24111
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H11").Value = 100373.5
$ws_CUL.Range("I11").Value = 318.5
$ws_CUL.Range("J11").Value = 250456
$ws_CUL.Range("K11").Value = 955.5
$ws_CUL.Range("L11").Value = 751368
$ws_CUL.Range("M11").Value = -815.5
$ws_CUL.Range("N11").Value = -751648
$ws_CUL.Range("H12").Value = 520
$ws_CUL.Range("J12").Value = 520
$ws_CUL.Range("L12").Value = 1560
$ws_CUL.Range("N12").Value = -1906
$ws_CUL.Range("H17").Value = 185.95
$ws_CUL.Range("I17").Value = 125
$ws_CUL.Range("J17").Value = 246.9
$ws_CUL.Range("K17").Value = 375
$ws_CUL.Range("L17").Value = 740.7
$ws_CUL.Range("M17").Value = -206
$ws_CUL.Range("N17").Value = -1078.7
$ws_CUL.Range("H68").Value = 1698.4286
$ws_CUL.Range("I68").Value = 1111.6
$ws_CUL.Range("J68").Value = 3165.5
$ws_CUL.Range("K68").Value = 3334.8
$ws_CUL.Range("L68").Value = 9496.5
$ws_CUL.Range("M68").Value = -2523.8
$ws_CUL.Range("N68").Value = -11118.5
$ws_CUL.Range("H71").Value = 1698.4286
$ws_CUL.Range("I71").Value = 1111.6
$ws_CUL.Range("J71").Value = 3165.5
$ws_CUL.Range("K71").Value = 10004.4
$ws_CUL.Range("L71").Value = 28489.5
$ws_CUL.Range("M71").Value = -5948.4
$ws_CUL.Range("N71").Value = -36601.5
$ws_CUL.Range("H93").Value = 5218.75
$ws_CUL.Range("J93").Value = 5218.75
$ws_CUL.Range("L93").Value = 15656.25
$ws_CUL.Range("N93").Value = -19400.25
$ws_CUL.Range("H114").Value = 18182784
$ws_CUL.Range("I114").Value = 22222848
$ws_CUL.Range("K114").Value = 66668544
$ws_CUL.Range("M114").Value = -66665290
$ws_CUL.Range("H117").Value = 4737
$ws_CUL.Range("I117").Value = 395.5
$ws_CUL.Range("J117").Value = 7105.091
$ws_CUL.Range("K117").Value = 1186.5
$ws_CUL.Range("L117").Value = 21315.273
$ws_CUL.Range("M117").Value = 2255.5
$ws_CUL.Range("N117").Value = -28199.273
$ws_CUL.Range("H121").Value = 33434984
$ws_CUL.Range("I121").Value = 83334264
$ws_CUL.Range("K121").Value = 250002792
$ws_CUL.Range("M121").Value = -250001482
$ws_CUL.Range("H128").Value = 316790.34
$ws_CUL.Range("I128").Value = 316790.34
$ws_CUL.Range("K128").Value = 950371.02
$ws_CUL.Range("M128").Value = -945391.02
$ws_CUL.Range("H140").Value = 4252.269
$ws_CUL.Range("I140").Value = 2012.5555
$ws_CUL.Range("K140").Value = 6037.666499999999
$ws_CUL.Range("M140").Value = -857.6664999999994
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H51").Value = 0
$ws_GSM.Range("J51").Value = 0
$ws_GSM.Range("L51").Value = 0
$ws_GSM.Range("N51").ClearContents()
$ws_GSM.Range("H70").Value = 16671558
$ws_GSM.Range("I70").Value = 22226954
$ws_GSM.Range("K70").Value = 22226954
$ws_GSM.Range("M70").Value = -22226684
$ws_GSM.Range("H73").Value = 16671558
$ws_GSM.Range("I73").Value = 22226954
$ws_GSM.Range("K73").Value = 22226954
$ws_GSM.Range("M73").Value = -22226018
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H16").Value = 724.125
$ws_LTW.Range("I16").Value = 724.125
$ws_LTW.Range("K16").Value = 724.125
$ws_LTW.Range("M16").Value = -554.125
$ws_LTW.Range("H122").Value = 2820.6667
$ws_LTW.Range("I122").Value = 2690.7778
$ws_LTW.Range("J122").Value = 3600
$ws_LTW.Range("K122").Value = 8072.3334
$ws_LTW.Range("L122").Value = 10800
$ws_LTW.Range("M122").Value = -5622.3334
$ws_LTW.Range("N122").Value = -15700
$ws_LTW.Range("H132").Value = 7840.4546
$ws_LTW.Range("I132").Value = 13435.111
$ws_LTW.Range("J132").Value = 3967.2307
$ws_LTW.Range("K132").Value = 40305.333
$ws_LTW.Range("L132").Value = 11901.6921
$ws_LTW.Range("M132").Value = -37775.333
$ws_LTW.Range("N132").Value = -16961.6921
$ws_LTW.Range("H136").Value = 3926.182
$ws_LTW.Range("I136").Value = 3673.25
$ws_LTW.Range("K136").Value = 11019.75
$ws_LTW.Range("M136").Value = -8469.75
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H8").Value = 1750
$ws_WVR.Range("I8").Value = 1000
$ws_WVR.Range("J8").Value = 2500
$ws_WVR.Range("K8").Value = 1000
$ws_WVR.Range("L8").Value = 2500
$ws_WVR.Range("M8").Value = -860
$ws_WVR.Range("N8").Value = -2780
$ws_WVR.Range("H22").Value = 9341.666999999999
$ws_WVR.Range("I22").Value = 9505
$ws_WVR.Range("K22").Value = 9505
$ws_WVR.Range("M22").Value = -9212
$ws_WVR.Range("H52").Value = 19999.375
$ws_WVR.Range("I52").Value = 19000
$ws_WVR.Range("J52").Value = 20332.5
$ws_WVR.Range("K52").Value = 19000
$ws_WVR.Range("L52").Value = 20332.5
$ws_WVR.Range("M52").Value = -18774
$ws_WVR.Range("N52").Value = -20784.5
$ws_WVR.Range("H113").Value = 526.73914
$ws_WVR.Range("I113").Value = 335.4
$ws_WVR.Range("K113").Value = 1006.2
$ws_WVR.Range("M113").Value = 1163.8
$ws_WVR.Range("H122").Value = 1555.2307
$ws_WVR.Range("J122").Value = 2255.923
$ws_WVR.Range("L122").Value = 6767.768999999999
$ws_WVR.Range("N122").Value = -11667.769
$ws_WVR.Range("H132").Value = 30999.666
$ws_WVR.Range("I132").Value = 26427.715
$ws_WVR.Range("J132").Value = 47001.5
$ws_WVR.Range("K132").Value = 79283.145
$ws_WVR.Range("L132").Value = 141004.5
$ws_WVR.Range("M132").Value = -76753.145
$ws_WVR.Range("N132").Value = -146064.5
